# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
# Cells that hold numeric-looking text (e.g. "1.004") must be written with a
# text NumberFormat first, otherwise Excel auto-converts them to real numbers;
# the format is reset back to Normal afterwards so no stray styling is left
# behind (matches the source workbook, where these are plain inline strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.925.69"
$ws.Range("E2").Value = "  -0.36%  "
# Row 3
$ws.Range("D3").Value = "1.856.57"
$ws.Range("E3").Value = "  -1.48%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5120"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.35%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3805"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.00%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.52%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.169"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.58%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.862.57"
$ws.Range("E13").Value = "  -1.30%  "
# Row 14
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001092"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.07%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.39%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06598"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "
# Row 20
$ws.Range("E20").Value = "  -1.69%  "
# Row 21
$ws.Range("E21").Value = "  +0.25%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.08%  "
# Row 23
$ws.Range("D23").Value = "27.961.30"
$ws.Range("E23").Value = "  -0.38%  "
# Row 24
$ws.Range("E24").Value = "  -3.62%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.13%  "
# Row 27
$ws.Range("D27").Value = "2.070.91"
$ws.Range("E27").Value = "  -1.42%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.13%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.10%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "
# Row 32
$ws.Range("E32").Value = "  -2.83%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.596"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "
# Row 35
$ws.Range("E35").Value = "  +1.54%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06497"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.24%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02401"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2159"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.200"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.79%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6426"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
# Row 41
$ws.Range("E41").Value = "  -4.60%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.07%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.853"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.36%  "
# Row 46
$ws.Range("E46").Value = "  -0.92%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.658"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.965"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.34%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.23%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.41%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.03%  "
